$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item('展览')

# F-column (想去人数) updates for existing rows
$fChanges = @{
  2 = 1629
  3 = 857
  4 = 261
  5 = 80
  6 = 1171
  7 = 787
  8 = 816
  9 = 1507
  10 = 299
  11 = 1050
  13 = 72
  14 = 198
  15 = 57
  16 = 499
  17 = 53
  18 = 38
  19 = 7
  22 = 572
  23 = 578
  24 = 45
  26 = 772
}
foreach ($r in $fChanges.Keys) { $ws.Range('F' + $r).Value = $fChanges[$r] }

# Insert new row 29 (event on 2024-04-13), shifting the old row 29 down to row 30
$ws.Rows.Item(29).Insert()
$ws.Range('A28:I28').Copy($ws.Range('A29:I29'))

$ws.Range('B29').NumberFormat = '@'
$ws.Range('A29').Value = 28
$ws.Range('B29').Value = '2024-04-13'
$ws.Range('C29').Value = '广州·Veni Vidi Vici动漫游戏嘉年华'
$ws.Range('D29').Value = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$ws.Range('E29').Value = '2024.04.13 10:00-04.13 17:00'
$ws.Range('F29').Value = 0
$ws.Range('G29').Value = '不可售'
$ws.Range('H29').Value = 'https://show.bilibili.com/platform/detail.html?id=81575'
$ws.Range('I29').Value = '//i1.hdslb.com/bfs/openplatform/202401/7ir7DZHt1706697841803.jpeg'

# Fix renumbered row 30 (was row 29 pre-insert): index + F value
$ws.Range('A30').Value = 29
$ws.Range('F30').Value = 374

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item('演出')
$fChanges = @{
  3 = 1023
  6 = 17
  8 = 69
  12 = 7
}
foreach ($r in $fChanges.Keys) { $ws.Range('F' + $r).Value = $fChanges[$r] }

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item('本地生活')
$fChanges = @{
  2 = 263
}
foreach ($r in $fChanges.Keys) { $ws.Range('F' + $r).Value = $fChanges[$r] }

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item('全部类型')

# F-column (想去人数) updates for existing (unshifted) rows
$fChanges = @{
  2 = 263
  3 = 1629
  5 = 857
  6 = 261
  7 = 1023
  8 = 80
  9 = 1171
  10 = 787
  11 = 816
  12 = 1507
  13 = 299
  14 = 1050
  16 = 72
  17 = 198
  18 = 57
  19 = 499
  20 = 53
  21 = 38
  23 = 7
  27 = 17
  30 = 572
  31 = 578
  32 = 45
  34 = 772
  36 = 69
}
foreach ($r in $fChanges.Keys) { $ws.Range('F' + $r).Value = $fChanges[$r] }

# Insert new row 41 (event on 2024-04-13), shifting rows 41-43 down to 42-44
$ws.Rows.Item(41).Insert()
$ws.Range('A40:I40').Copy($ws.Range('A41:I41'))

$ws.Range('B41').NumberFormat = '@'
$ws.Range('A41').Value = 40
$ws.Range('B41').Value = '2024-04-13'
$ws.Range('C41').Value = '广州·Veni Vidi Vici动漫游戏嘉年华'
$ws.Range('D41').Value = '广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心'
$ws.Range('E41').Value = '2024.04.13 10:00-04.13 17:00'
$ws.Range('F41').Value = 0
$ws.Range('G41').Value = '不可售'
$ws.Range('H41').Value = 'https://show.bilibili.com/platform/detail.html?id=81575'
$ws.Range('I41').Value = '//i1.hdslb.com/bfs/openplatform/202401/7ir7DZHt1706697841803.jpeg'

# Fix renumbered rows 42-44 (index column) and F-values shifted with them
$ws.Range('A42').Value = 41
$ws.Range('A43').Value = 42
$ws.Range('A44').Value = 43
$ws.Range('F43').Value = 374
$ws.Range('F44').Value = 7

